$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @(3,22376,3565,3761),
    @(4,11178,1536,1395),
    @(5,39716,3553,3791),
    @(6,643,358,84),
    @(7,25106,4583,3892),
    @(8,3181,730,700),
    @(9,3107,582,395),
    @(10,1331,233,160),
    @(11,94,126,22),
    @(12,3,0,0),
    @(13,555,136,184),
    @(14,1605,720,459),
    @(15,2873,1084,478),
    @(16,1766,685,251),
    @(17,876,442,80),
    @(18,8803,1613,1760),
    @(19,1021,356,342),
    @(20,10118,1230,1757),
    @(21,121,210,7),
    @(22,9031,1245,1761),
    @(23,682,242,79),
    @(24,9980,1648,2126),
    @(25,43686,4148,5459),
    @(26,2835,1009,577),
    @(27,0,0,0),
    @(28,2869,571,737),
    @(29,733,273,162),
    @(30,7513,1519,1322),
    @(31,283,106,144),
    @(32,1375,967,266),
    @(33,8430,1939,1507),
    @(34,5206,1682,1276),
    @(35,3223,378,891),
    @(36,28978,3254,3113),
    @(37,4244,1663,697),
    @(38,12804,1169,1584),
    @(39,490,505,154),
    @(40,1134,284,439),
    @(41,2112,253,90),
    @(42,7913,450,246),
    @(43,219,78,65),
    @(44,534,38,43),
    @(45,1045,14,2),
    @(46,1650,533,224),
    @(47,6143,1906,1140),
    @(48,16473,2074,2579),
    @(49,6786,1983,653),
    @(50,6205,607,886),
    @(51,15921,1678,2548),
    @(52,2568,312,708),
    @(53,7947,1680,1503),
    @(54,854,645,356),
    @(55,1202,784,86),
    @(56,1365,398,417),
    @(57,6799,2625,1372),
    @(58,11823,833,446),
    @(59,356906,57777,53358)
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
}
